# Atualizado por script em 26-11-2023 20:30
# Adds 3 new match rows (124-126) to the Paraguay Primera Division 2023 sheet,
# mirroring the formatting of the last existing data row (123).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 124 (Indice 123) ---------------------------------------------
$ws.Range("A123:V123").Copy($ws.Range("A124:V124"))
$ws.Range("A124").Value = 123
$ws.Range("E124").Value = 45254.97916666666
$ws.Range("F124").Value = "Ameliano"
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = "General Caballero JLM"
$ws.Range("I124").Value = 2
$ws.Range("J124").Value = 1.72
$ws.Range("K124").Value = "19/11/2023 01:12"
$ws.Range("L124").Value = 2.15
$ws.Range("M124").Value = "24/11/2023 23:21"
$ws.Range("N124").Value = 3.79
$ws.Range("O124").Value = "19/11/2023 01:12"
$ws.Range("P124").Value = 3.32
$ws.Range("Q124").Value = "24/11/2023 23:23"
$ws.Range("R124").Value = 4.93
$ws.Range("S124").Value = "19/11/2023 01:12"
$ws.Range("T124").Value = 3.72
$ws.Range("U124").Value = "24/11/2023 23:21"
$ws.Range("V124").Value = "https://www.betexplorer.com/football/paraguay/primera-division/sportivo-ameliano-general-caballero-jlm/j95nslHf/"

# --- Row 125 (Indice 124) ---------------------------------------------
$ws.Range("A123:V123").Copy($ws.Range("A125:V125"))
$ws.Range("A125").Value = 124
$ws.Range("E125").Value = 45255.97916666666
$ws.Range("F125").Value = "Guarani"
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = "Sportivo Trinidense"
$ws.Range("I125").Value = 2
$ws.Range("J125").Value = 2.23
$ws.Range("K125").Value = "19/11/2023 01:13"
$ws.Range("L125").Value = 2.01
$ws.Range("M125").Value = "25/11/2023 23:22"
$ws.Range("N125").Value = 3.55
$ws.Range("O125").Value = "19/11/2023 01:13"
$ws.Range("P125").Value = 3.63
$ws.Range("Q125").Value = "25/11/2023 23:22"
$ws.Range("R125").Value = 3.02
$ws.Range("S125").Value = "19/11/2023 01:13"
$ws.Range("T125").Value = 3.81
$ws.Range("U125").Value = "25/11/2023 23:22"
$ws.Range("V125").Value = "https://www.betexplorer.com/football/paraguay/primera-division/guarani-sportivo-trinidense/Ua4jt8W0/"

# --- Row 126 (Indice 125) ---------------------------------------------
$ws.Range("A123:V123").Copy($ws.Range("A126:V126"))
$ws.Range("A126").Value = 125
$ws.Range("E126").Value = 45255.97916666666
$ws.Range("F126").Value = "Olimpia Asuncion"
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = "Libertad Asuncion"
$ws.Range("I126").Value = 3
$ws.Range("J126").Value = 2.99
$ws.Range("K126").Value = "19/11/2023 01:12"
$ws.Range("L126").Value = 3.2
$ws.Range("M126").Value = "25/11/2023 23:29"
$ws.Range("N126").Value = 3.3
$ws.Range("O126").Value = "19/11/2023 01:12"
$ws.Range("P126").Value = 3.18
$ws.Range("Q126").Value = "25/11/2023 23:29"
$ws.Range("R126").Value = 2.47
$ws.Range("S126").Value = "19/11/2023 01:12"
$ws.Range("T126").Value = 2.46
$ws.Range("U126").Value = "25/11/2023 23:29"
$ws.Range("V126").Value = "https://www.betexplorer.com/football/paraguay/primera-division/olimpia-asuncion-libertad-asuncion/COgeuSo7/"

Write-Output "Added rows 124-126"
